# Apply scraped schedule update for Linea 141 workbook
# (horarios-141-2026-01-15.xlsx) - commit: Horarios actualizados Linea 141 - 389

$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 14:44:25"
$ws.Cells.Item(3, 1).Value = "Total filas: 283"

# Row 128
$ws.Cells.Item(128, 1).Value = "08:36:20"
$ws.Cells.Item(128, 2).Value = "10:22"
$ws.Cells.Item(128, 3).Value = "17_ROMERO"
$ws.Cells.Item(128, 4).Value = 106
$ws.Cells.Item(128, 5).Value = "LP1912"
# Row 129
$ws.Cells.Item(129, 1).Value = "09:25:56"
$ws.Cells.Item(129, 2).Value = "10:22"
$ws.Cells.Item(129, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(129, 4).Value = 57
$ws.Cells.Item(129, 5).Value = "LP1912"
# Row 146
$ws.Cells.Item(146, 1).Value = "10:12:35"
$ws.Cells.Item(146, 2).Value = "10:56"
$ws.Cells.Item(146, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(146, 4).Value = 44
$ws.Cells.Item(146, 5).Value = "LP1912"
# Row 147
$ws.Cells.Item(147, 1).Value = "10:52:48"
$ws.Cells.Item(147, 2).Value = "10:56"
$ws.Cells.Item(147, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(147, 4).Value = 4
$ws.Cells.Item(147, 5).Value = "LP1912"
# Row 197
$ws.Cells.Item(197, 1).Value = "10:52:48"
$ws.Cells.Item(197, 2).Value = "12:21"
$ws.Cells.Item(197, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(197, 4).Value = 89
$ws.Cells.Item(197, 5).Value = "LP1912"
# Row 198
$ws.Cells.Item(198, 1).Value = "11:17:08"
$ws.Cells.Item(198, 2).Value = "12:21"
$ws.Cells.Item(198, 3).Value = "215A_EL PATO"
$ws.Cells.Item(198, 4).Value = 64
$ws.Cells.Item(198, 5).Value = "LP1912"
# Row 199
$ws.Cells.Item(199, 1).Value = "11:17:08"
$ws.Cells.Item(199, 2).Value = "12:21"
$ws.Cells.Item(199, 3).Value = "14_ABASTO"
$ws.Cells.Item(199, 4).Value = 64
$ws.Cells.Item(199, 5).Value = "LP1912"
# Row 243
$ws.Cells.Item(243, 1).Value = "12:27:08"
$ws.Cells.Item(243, 2).Value = "14:04"
$ws.Cells.Item(243, 3).Value = "17_ROMERO"
$ws.Cells.Item(243, 4).Value = 97
$ws.Cells.Item(243, 5).Value = "LP1912"
# Row 244
$ws.Cells.Item(244, 1).Value = "14:00:52"
$ws.Cells.Item(244, 2).Value = "14:04"
$ws.Cells.Item(244, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(244, 4).Value = 4
$ws.Cells.Item(244, 5).Value = "LP1912"
# Row 255
$ws.Cells.Item(255, 1).Value = "14:44:25"
$ws.Cells.Item(255, 2).Value = "14:45"
$ws.Cells.Item(255, 3).Value = "15_ABASTO"
$ws.Cells.Item(255, 4).Value = 1
$ws.Cells.Item(255, 5).Value = "LP1912"
# Row 256
$ws.Cells.Item(256, 1).Value = "14:00:52"
$ws.Cells.Item(256, 2).Value = "14:56"
$ws.Cells.Item(256, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(256, 4).Value = 56
$ws.Cells.Item(256, 5).Value = "LP1912"
# Row 257
$ws.Cells.Item(257, 1).Value = "13:23:09"
$ws.Cells.Item(257, 2).Value = "14:57"
$ws.Cells.Item(257, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(257, 4).Value = 94
$ws.Cells.Item(257, 5).Value = "LP1912"
# Row 258
$ws.Cells.Item(258, 1).Value = "13:23:09"
$ws.Cells.Item(258, 2).Value = "14:58"
$ws.Cells.Item(258, 3).Value = "215B_EL PATO"
$ws.Cells.Item(258, 4).Value = 95
$ws.Cells.Item(258, 5).Value = "LP1912"
# Row 259
$ws.Cells.Item(259, 1).Value = "13:23:09"
$ws.Cells.Item(259, 2).Value = "15:00"
$ws.Cells.Item(259, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(259, 4).Value = 97
$ws.Cells.Item(259, 5).Value = "LP1912"
# Row 260
$ws.Cells.Item(260, 1).Value = "14:00:52"
$ws.Cells.Item(260, 2).Value = "15:04"
$ws.Cells.Item(260, 3).Value = "10_OLMOS"
$ws.Cells.Item(260, 4).Value = 64
$ws.Cells.Item(260, 5).Value = "LP1912"
# Row 261
$ws.Cells.Item(261, 1).Value = "14:44:25"
$ws.Cells.Item(261, 2).Value = "15:04"
$ws.Cells.Item(261, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(261, 4).Value = 20
$ws.Cells.Item(261, 5).Value = "LP1912"
# Row 262
$ws.Cells.Item(262, 1).Value = "13:23:09"
$ws.Cells.Item(262, 2).Value = "15:05"
$ws.Cells.Item(262, 3).Value = "10_OLMOS"
$ws.Cells.Item(262, 4).Value = 102
$ws.Cells.Item(262, 5).Value = "LP1912"
# Row 263
$ws.Cells.Item(263, 1).Value = "14:44:25"
$ws.Cells.Item(263, 2).Value = "15:06"
$ws.Cells.Item(263, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(263, 4).Value = 22
$ws.Cells.Item(263, 5).Value = "LP1912"
# Row 264
$ws.Cells.Item(264, 1).Value = "14:00:52"
$ws.Cells.Item(264, 2).Value = "15:10"
$ws.Cells.Item(264, 3).Value = "17_ROMERO"
$ws.Cells.Item(264, 4).Value = 70
$ws.Cells.Item(264, 5).Value = "LP1912"
# Row 265
$ws.Cells.Item(265, 1).Value = "14:00:52"
$ws.Cells.Item(265, 2).Value = "15:13"
$ws.Cells.Item(265, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(265, 4).Value = 73
$ws.Cells.Item(265, 5).Value = "LP1912"
# Row 266
$ws.Cells.Item(266, 1).Value = "13:23:09"
$ws.Cells.Item(266, 2).Value = "15:14"
$ws.Cells.Item(266, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(266, 4).Value = 111
$ws.Cells.Item(266, 5).Value = "LP1912"
# Row 267
$ws.Cells.Item(267, 1).Value = "14:44:25"
$ws.Cells.Item(267, 2).Value = "15:16"
$ws.Cells.Item(267, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(267, 4).Value = 32
$ws.Cells.Item(267, 5).Value = "LP1912"
# Row 268
$ws.Cells.Item(268, 1).Value = "14:00:52"
$ws.Cells.Item(268, 2).Value = "15:20"
$ws.Cells.Item(268, 3).Value = "15_ABASTO"
$ws.Cells.Item(268, 4).Value = 80
$ws.Cells.Item(268, 5).Value = "LP1912"
# Row 269
$ws.Cells.Item(269, 1).Value = "13:23:09"
$ws.Cells.Item(269, 2).Value = "15:21"
$ws.Cells.Item(269, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(269, 4).Value = 118
$ws.Cells.Item(269, 5).Value = "LP1912"
# Row 270
$ws.Cells.Item(270, 1).Value = "14:00:52"
$ws.Cells.Item(270, 2).Value = "15:25"
$ws.Cells.Item(270, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(270, 4).Value = 85
$ws.Cells.Item(270, 5).Value = "LP1912"
# Row 271
$ws.Cells.Item(271, 1).Value = "14:00:52"
$ws.Cells.Item(271, 2).Value = "15:32"
$ws.Cells.Item(271, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(271, 4).Value = 92
$ws.Cells.Item(271, 5).Value = "LP1912"
# Row 272
$ws.Cells.Item(272, 1).Value = "14:00:52"
$ws.Cells.Item(272, 2).Value = "15:35"
$ws.Cells.Item(272, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(272, 4).Value = 95
$ws.Cells.Item(272, 5).Value = "LP1912"
# Row 273
$ws.Cells.Item(273, 1).Value = "14:00:52"
$ws.Cells.Item(273, 2).Value = "15:36"
$ws.Cells.Item(273, 3).Value = "10_OLMOS"
$ws.Cells.Item(273, 4).Value = 96
$ws.Cells.Item(273, 5).Value = "LP1912"
# Row 274
$ws.Cells.Item(274, 1).Value = "14:44:25"
$ws.Cells.Item(274, 2).Value = "15:37"
$ws.Cells.Item(274, 3).Value = "10_OLMOS"
$ws.Cells.Item(274, 4).Value = 53
$ws.Cells.Item(274, 5).Value = "LP1912"
# Row 275
$ws.Cells.Item(275, 1).Value = "14:00:52"
$ws.Cells.Item(275, 2).Value = "15:38"
$ws.Cells.Item(275, 3).Value = "215A_EL PATO"
$ws.Cells.Item(275, 4).Value = 98
$ws.Cells.Item(275, 5).Value = "LP1912"
# Row 276
$ws.Cells.Item(276, 1).Value = "14:44:25"
$ws.Cells.Item(276, 2).Value = "15:45"
$ws.Cells.Item(276, 3).Value = "14_ABASTO"
$ws.Cells.Item(276, 4).Value = 61
$ws.Cells.Item(276, 5).Value = "LP1912"
# Row 277
$ws.Cells.Item(277, 1).Value = "14:00:52"
$ws.Cells.Item(277, 2).Value = "15:46"
$ws.Cells.Item(277, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(277, 4).Value = 106
$ws.Cells.Item(277, 5).Value = "LP1912"
# Row 278
$ws.Cells.Item(278, 1).Value = "14:44:25"
$ws.Cells.Item(278, 2).Value = "15:47"
$ws.Cells.Item(278, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(278, 4).Value = 63
$ws.Cells.Item(278, 5).Value = "LP1912"
# Row 279
$ws.Cells.Item(279, 1).Value = "14:00:52"
$ws.Cells.Item(279, 2).Value = "15:53"
$ws.Cells.Item(279, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(279, 4).Value = 113
$ws.Cells.Item(279, 5).Value = "LP1912"
# Row 280
$ws.Cells.Item(280, 1).Value = "14:44:25"
$ws.Cells.Item(280, 2).Value = "15:56"
$ws.Cells.Item(280, 3).Value = "17_ROMERO"
$ws.Cells.Item(280, 4).Value = 72
$ws.Cells.Item(280, 5).Value = "LP1912"
# Row 281
$ws.Cells.Item(281, 1).Value = "14:00:52"
$ws.Cells.Item(281, 2).Value = "15:56"
$ws.Cells.Item(281, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(281, 4).Value = 116
$ws.Cells.Item(281, 5).Value = "LP1912"
# Row 282
$ws.Cells.Item(282, 1).Value = "14:44:25"
$ws.Cells.Item(282, 2).Value = "16:01"
$ws.Cells.Item(282, 3).Value = "10_OLMOS"
$ws.Cells.Item(282, 4).Value = 77
$ws.Cells.Item(282, 5).Value = "LP1912"
# Row 283
$ws.Cells.Item(283, 1).Value = "14:44:25"
$ws.Cells.Item(283, 2).Value = "16:02"
$ws.Cells.Item(283, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(283, 4).Value = 78
$ws.Cells.Item(283, 5).Value = "LP1912"
# Row 284
$ws.Cells.Item(284, 1).Value = "14:44:25"
$ws.Cells.Item(284, 2).Value = "16:15"
$ws.Cells.Item(284, 3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(284, 4).Value = 91
$ws.Cells.Item(284, 5).Value = "LP1912"
# Row 285
$ws.Cells.Item(285, 1).Value = "14:44:25"
$ws.Cells.Item(285, 2).Value = "16:20"
$ws.Cells.Item(285, 3).Value = "215C_EL PATO"
$ws.Cells.Item(285, 4).Value = 96
$ws.Cells.Item(285, 5).Value = "LP1912"
# Row 286
$ws.Cells.Item(286, 1).Value = "14:44:25"
$ws.Cells.Item(286, 2).Value = "16:21"
$ws.Cells.Item(286, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(286, 4).Value = 97
$ws.Cells.Item(286, 5).Value = "LP1912"
# Row 287
$ws.Cells.Item(287, 1).Value = "14:44:25"
$ws.Cells.Item(287, 2).Value = "16:42"
$ws.Cells.Item(287, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(287, 4).Value = 118
$ws.Cells.Item(287, 5).Value = "LP1912"
# Row 288
$ws.Cells.Item(288, 1).Value = "14:44:25"
$ws.Cells.Item(288, 2).Value = "16:43"
$ws.Cells.Item(288, 3).Value = "225_GOMEZ"
$ws.Cells.Item(288, 4).Value = 119
$ws.Cells.Item(288, 5).Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 14:44:25"
$ws.Cells.Item(3, 1).Value = "Total filas: 29"

# Row 34
$ws.Cells.Item(34, 1).Value = "14:44:25"
$ws.Cells.Item(34, 2).Value = "16:20"
$ws.Cells.Item(34, 3).Value = "215C_EL PATO"
$ws.Cells.Item(34, 4).Value = 96
$ws.Cells.Item(34, 5).Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 14:44:25"
$ws.Cells.Item(3, 1).Value = "Total filas: 39"

# Row 40
$ws.Cells.Item(40, 1).Value = "14:44:25"
$ws.Cells.Item(40, 2).Value = "14:46"
$ws.Cells.Item(40, 3).Value = "215D_LA PLATA"
$ws.Cells.Item(40, 4).Value = 2
$ws.Cells.Item(40, 5).Value = "L6203"
# Row 41
$ws.Cells.Item(41, 1).Value = "14:00:52"
$ws.Cells.Item(41, 2).Value = "14:52"
$ws.Cells.Item(41, 3).Value = "215D_LA PLATA"
$ws.Cells.Item(41, 4).Value = 52
$ws.Cells.Item(41, 5).Value = "L6203"
# Row 42
$ws.Cells.Item(42, 1).Value = "12:54:06"
$ws.Cells.Item(42, 2).Value = "14:53"
$ws.Cells.Item(42, 3).Value = "215D_LA PLATA"
$ws.Cells.Item(42, 4).Value = 119
$ws.Cells.Item(42, 5).Value = "L6203"
# Row 43
$ws.Cells.Item(43, 1).Value = "14:00:52"
$ws.Cells.Item(43, 2).Value = "15:34"
$ws.Cells.Item(43, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(43, 4).Value = 94
$ws.Cells.Item(43, 5).Value = "L6173"
# Row 44
$ws.Cells.Item(44, 1).Value = "14:44:25"
$ws.Cells.Item(44, 2).Value = "16:14"
$ws.Cells.Item(44, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(44, 4).Value = 90
$ws.Cells.Item(44, 5).Value = "L6203"

